$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.190.28"
$ws.Range("E2").Value = "  -0.76%  "

$ws.Range("D3").Value = "1.857.82"
$ws.Range("E3").Value = "  -0.81%  "

$ws.Range("D4").Value = "'1.0000"
$ws.Range("E4").Value = "  -0.11%  "

$ws.Range("D5").Value = "'0.7132"
$ws.Range("E5").Value = "  -0.47%  "

$ws.Range("D6").Value = "'240.29"
$ws.Range("E6").Value = "  +0.50%  "

$ws.Range("D7").Value = "'0.9999"

$ws.Range("D8").Value = "'0.07751"
$ws.Range("E8").Value = "  -0.77%  "

$ws.Range("E9").Value = "  +0.13%  "

$ws.Range("D10").Value = "'25.07"
$ws.Range("E10").Value = "  -0.88%  "

$ws.Range("D11").Value = "'0.08251"

$ws.Range("B12").Value = "Polkadot"
$ws.Range("C12").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D12").Value = "'5.230"
$ws.Range("E12").Value = "  -0.19%  "

$ws.Range("B13").Value = "Polygon"
$ws.Range("C13").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D13").Value = "'0.7156"
$ws.Range("E13").Value = "  -1.00%  "

$ws.Range("B14").Value = "WrappedEther"
$ws.Range("C14").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D14").Value = "1.841.30"
$ws.Range("E14").Value = "  -1.86%  "

$ws.Range("D15").Value = "'90.23"
$ws.Range("E15").Value = "  -0.34%  "

$ws.Range("D16").Value = "29.174.59"
$ws.Range("E16").Value = "  -1.04%  "

$ws.Range("D17").Value = "'5.860"
$ws.Range("E17").Value = "  +0.15%  "

$ws.Range("D18").Value = "'243.57"
$ws.Range("E18").Value = "  +0.69%  "

$ws.Range("D19").Value = "'0.000007790"
$ws.Range("E19").Value = "  -0.97%  "

$ws.Range("D20").Value = "'13.15"
$ws.Range("E20").Value = "  -0.96%  "

$ws.Range("D21").Value = "2.106.07"
$ws.Range("E21").Value = "  -1.31%  "

$ws.Range("D22").Value = "'0.9997"
$ws.Range("E22").Value = "  -0.03%  "

$ws.Range("D23").Value = "'7.998"
$ws.Range("E23").Value = "  +3.13%  "

$ws.Range("E24").Value = "  -0.13%  "

$ws.Range("D25").Value = "'0.1596"
$ws.Range("E25").Value = "  +2.88%  "

$ws.Range("D26").Value = "'162.39"
$ws.Range("E26").Value = "  -0.45%  "

$ws.Range("D27").Value = "'8.906"
$ws.Range("E27").Value = "  -1.01%  "

$ws.Range("D28").Value = "'18.31"
$ws.Range("E28").Value = "  -0.12%  "

$ws.Range("E29").Value = "  +0.92%  "

$ws.Range("E30").Value = "  -3.11%  "

$ws.Range("D31").Value = "'4.401"
$ws.Range("E31").Value = "  +1.67%  "

$ws.Range("D32").Value = "'4.202"
$ws.Range("E32").Value = "  +2.76%  "

$ws.Range("D33").Value = "'0.05184"
$ws.Range("E33").Value = "  -1.21%  "

$ws.Range("E34").Value = "  -1.14%  "

$ws.Range("E35").Value = "  -2.29%  "

$ws.Range("D36").Value = "'0.7262"
$ws.Range("E36").Value = "  +1.19%  "

$ws.Range("D37").Value = "'2.676"
$ws.Range("E37").Value = "  -0.04%  "

$ws.Range("E38").Value = "  -0.54%  "

$ws.Range("D39").Value = "'2.686"
$ws.Range("E39").Value = "  -1.46%  "

$ws.Range("D40").Value = "1.162.73"
$ws.Range("E40").Value = "  -1.69%  "

$ws.Range("D41").Value = "'0.9050"
$ws.Range("E41").Value = "  -0.32%  "

$ws.Range("D42").Value = "'6.150"
$ws.Range("E42").Value = "  +2.29%  "

$ws.Range("D43").Value = "'72.18"
$ws.Range("E43").Value = "  -0.03%  "

$ws.Range("D44").Value = "'0.9996"
$ws.Range("E44").Value = "  -0.12%  "

$ws.Range("D45").Value = "'101.55"
$ws.Range("E45").Value = "  -0.77%  "

$ws.Range("D46").Value = "2.002.02"
$ws.Range("E46").Value = "  -1.09%  "

$ws.Range("D47").Value = "'0.5215"
$ws.Range("E47").Value = "  -2.89%  "

$ws.Range("D48").Value = "'1.764"
$ws.Range("E48").Value = "  +0.01%  "

$ws.Range("E49").Value = "  +1.03%  "

$ws.Range("D50").Value = "'9.308"
$ws.Range("E50").Value = "  +1.62%  "

$ws.Range("D51").Value = "'2.853"
$ws.Range("E51").Value = "  +1.27%  "
